$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.287.82"
$ws.Range("E2").Value = "  -1.30%  "

$ws.Range("D3").Value = "2.999.35"
$ws.Range("E3").Value = "  -1.74%  "

$ws.Range("D4").Value = "'" + "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'" + "588.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").Value = "'" + "145.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.69%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'" + "0.525"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.08%  "

$ws.Range("D9").Value = "2.995.70"
$ws.Range("E9").Value = "  -1.88%  "

$ws.Range("E10").Value = "  -4.37%  "

$ws.Range("D11").Value = "'" + "5.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.99%  "

$ws.Range("D12").Value = "'" + "0.467"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.33%  "

$ws.Range("D13").Value = "'" + "0.0000228"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.44%  "

$ws.Range("D14").Value = "'" + "34.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.84%  "

$ws.Range("E15").Value = "  +1.93%  "

$ws.Range("D16").Value = "3.499.47"
$ws.Range("E16").Value = "  -1.57%  "

$ws.Range("E17").Value = "  -0.96%  "

$ws.Range("D18").Value = "62.251.51"
$ws.Range("E18").Value = "  -1.33%  "

$ws.Range("D19").Value = "2.997.38"
$ws.Range("E19").Value = "  -1.79%  "

$ws.Range("D20").Value = "'" + "456.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.56%  "

$ws.Range("D21").Value = "'" + "14.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.40%  "

$ws.Range("D22").Value = "'" + "0.687"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.54%  "

$ws.Range("D23").Value = "'" + "7.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.59%  "

$ws.Range("D24").Value = "'" + "82.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("E25").Value = "  -9.29%  "

$ws.Range("D26").Value = "'" + "12.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.27%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("E28").Value = "  -8.19%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  -1.03%  "

$ws.Range("D31").Value = "'" + "6.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.11%  "

$ws.Range("E32").Value = "  -4.69%  "

$ws.Range("D33").Value = "'" + "27.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").Value = "'" + "0.108"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.00%  "

$ws.Range("D35").Value = "0.0₃0801"
$ws.Range("E35").Value = "  -1.94%  "

$ws.Range("E36").Value = "  -3.26%  "

$ws.Range("D37").Value = "'" + "5.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.62%  "

$ws.Range("E38").Value = "  -5.05%  "

$ws.Range("D39").Value = "'" + "9.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("D40").Value = "'" + "50.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("E41").Value = "  +7.21%  "

$ws.Range("D42").Value = "'" + "2.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.73%  "

$ws.Range("D43").Value = "'" + "392.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.59%  "

$ws.Range("E44").Value = "  -1.07%  "

$ws.Range("D45").Value = "'" + "0.266"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.29%  "

$ws.Range("D46").Value = "2.730.87"
$ws.Range("E46").Value = "  -3.43%  "

$ws.Range("D47").Value = "'" + "37.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.35%  "

$ws.Range("D48").Value = "'" + "129.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.27%  "

$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("E51").Value = "  -0.61%  "
